$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix trailing-space typos in existing shared strings (spell-check pass) ---
$ws.Range("D5").Value = "A Systematic Review of Cross- vs. Within-`nCompany Cost Estimation Studies"
$ws.Range("D6").Value = "Forecasting of Software De`nvelopment Work Effort: `nEvidence on Expert Judgment and Formal Models"
$ws.Range("E8").Value = "Relatório técnico sobre SLR em Engenharia de Software que propõe um conjunto de diretrizes"

# --- New row 9: entry #8 (08-BTH2012Yasin.pdf) ---
$ldquo = [char]0x201C
$rdquo = [char]0x201D
$ws.Range("C9").Value = "08-BTH2012Yasin.pdf"
$ws.Range("D9").Value = "On the quality of grey literature and its use`nin information synthesis`nduring systematic literature reviews"
$ws.Range("E9").Value = "Relatório de uma tese de mestrado sobre a utilização de " + $ldquo + "grey literature" + $rdquo + " em revisão sistemática"

# Row 9 grows to three lines of wrapped text, same as the other 3-line rows (row 6)
$ws.Rows(9).RowHeight = 35.05

# --- Selection moves to E10 ---
$ws.Range("E10").Select()

# --- tabColor normalizes from 00FFFFFF to FFFFFFFF (alpha channel) ---
$ws.Tab.Color = 16777215
